$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update header text for Start date / End date columns (F1 / G1)
$ws.Range("F1").Value = "Start date (m/d/y) or (d-m-y)"
$ws.Range("G1").Value = "End date (m/d/y) or (d-m-y)"

# Shift the End date column values forward by 239 days (2021-12-10 15:54 -> 2022-08-06 15:54)
for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.Value2 = $cell.Value2 + 239
}

# Widen column G slightly to fit the new header text
$ws.Columns.Item(7).ColumnWidth = 25.65

# Update the active selection / scroll position as saved by the author
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G5").Select()
